$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.177.47'
$ws.Range("E2").Value = '  +0.28%  '

# Row 3
$ws.Range("D3").Value = '2.620.38'
$ws.Range("E3").Value = '  -1.67%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.34'
$ws.Range("E5").Value = '  -0.18%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.68'
$ws.Range("E6").Value = '  +2.49%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.15%  '

# Row 8
$ws.Range("E8").Value = '  -2.18%  '

# Row 9
$ws.Range("D9").Value = '2.621.19'
$ws.Range("E9").Value = '  -1.56%  '

# Row 10
$ws.Range("E10").Value = '  +0.12%  '

# Row 11
$ws.Range("E11").Value = '  +1.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.361'
$ws.Range("E12").Value = '  +1.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.23'
$ws.Range("E13").Value = '  +0.83%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.59'
$ws.Range("E14").Value = '  -0.28%  '

# Row 15
$ws.Range("D15").Value = '3.107.72'
$ws.Range("E15").Value = '  -1.80%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("E16").Value = '  -0.10%  '

# Row 17
$ws.Range("D17").Value = '66.820.47'
$ws.Range("E17").Value = '  -0.46%  '

# Row 18
$ws.Range("D18").Value = '2.612.12'
$ws.Range("E18").Value = '  -2.34%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.12'
$ws.Range("E19").Value = '  +4.65%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.08'
$ws.Range("E20").Value = '  +7.97%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.41'
$ws.Range("E21").Value = '  -1.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  -1.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.66'
$ws.Range("E23").Value = '  -2.28%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.59'
$ws.Range("E24").Value = '  +6.41%  '

# Row 25
$ws.Range("E25").Value = '  -0.09%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.04'
$ws.Range("E26").Value = '  -1.80%  '

# Row 27
$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.90'
$ws.Range("E27").Value = '  -5.41%  '

# Row 28
$ws.Range("D28").Value = '2.758.39'
$ws.Range("E28").Value = '  -1.57%  '

# Row 29
$ws.Range("E29").Value = '  -0.38%  '

# Row 30
$ws.Range("E30").Value = '  -1.09%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '548.67'
$ws.Range("E31").Value = '  -0.27%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.93'
$ws.Range("E32").Value = '  -0.18%  '

# Row 33
$ws.Range("E33").Value = '  -1.88%  '

# Row 34
$ws.Range("E34").Value = '  -0.71%  '

# Row 35
$ws.Range("E35").Value = '  +5.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("E37").Value = '  -3.94%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.74'
$ws.Range("E38").Value = '  +0.74%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.00'
$ws.Range("E39").Value = '  -2.07%  '

# Row 40
$ws.Range("E40").Value = '  -1.48%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.17'
$ws.Range("E41").Value = '  -1.35%  '

# Row 42
$ws.Range("E42").Value = '  -1.57%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.90'
$ws.Range("E43").Value = '  -0.06%  '

# Row 44
$ws.Range("E44").Value = '  -0.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.21'
$ws.Range("E45").Value = '  -0.10%  '

# Row 46
$ws.Range("E46").Value = '  -4.54%  '

# Row 47
$ws.Range("D47").Value = '0.0₆0298'
$ws.Range("E47").Value = '  +0.28%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '151.34'
$ws.Range("E48").Value = '  -0.70%  '

# Row 49
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.577'
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.79'
$ws.Range("E50").Value = '  -0.68%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("E51").Value = '  -1.21%  '
